$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = '@'
$c.Value = '59.101.48'
$c.Style = 'Normal'
$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = '@'
$c.Value = '  +2.95%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = '@'
$c.Value = '2.502.38'
$c.Style = 'Normal'
$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = '@'
$c.Value = '  +3.17%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.34%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value = '534.79'
$c.Style = 'Normal'
$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = '@'
$c.Value = '  +4.30%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = '@'
$c.Value = '136.12'
$c.Style = 'Normal'
$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = '@'
$c.Value = '  +5.25%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.02%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = '@'
$c.Value = '0.568'
$c.Style = 'Normal'
$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = '@'
$c.Value = '  +3.66%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = '@'
$c.Value = '2.533.85'
$c.Style = 'Normal'
$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = '@'
$c.Value = '  +4.03%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = '@'
$c.Value = '  +5.26%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.93%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = '@'
$c.Value = '5.24'
$c.Style = 'Normal'
$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = '@'
$c.Value = '  +1.37%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = '@'
$c.Value = '2.962.74'
$c.Style = 'Normal'
$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = '@'
$c.Value = '  +3.75%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = '@'
$c.Value = '59.099.41'
$c.Style = 'Normal'
$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = '@'
$c.Value = '  +3.18%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = '@'
$c.Value = '22.60'
$c.Style = 'Normal'
$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = '@'
$c.Value = '  +3.99%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = '@'
$c.Value = '  +4.39%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = '@'
$c.Value = '2.530.48'
$c.Style = 'Normal'
$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = '@'
$c.Value = '  +4.06%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = '@'
$c.Value = '10.79'
$c.Style = 'Normal'
$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = '@'
$c.Value = '  +3.52%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = '@'
$c.Value = '  +3.85%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = '@'
$c.Value = '324.40'
$c.Style = 'Normal'
$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = '@'
$c.Value = '  +3.25%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = '@'
$c.Value = '  +8.84%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = '@'
$c.Value = '0.996'
$c.Style = 'Normal'
$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.34%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = '@'
$c.Value = '65.75'
$c.Style = 'Normal'
$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = '@'
$c.Value = '  +3.47%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = '@'
$c.Value = '  +1.36%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.35%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = '@'
$c.Value = '  +1.93%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(28, 5)
$c.NumberFormat = '@'
$c.Value = '  +5.81%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = '@'
$c.Value = '0.0₃0769'
$c.Style = 'Normal'
$c = $ws.Cells.Item(29, 5)
$c.NumberFormat = '@'
$c.Value = '  +6.89%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = '@'
$c.Value = '173.78'
$c.Style = 'Normal'
$c = $ws.Cells.Item(30, 5)
$c.NumberFormat = '@'
$c.Value = '  +3.00%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = '@'
$c.Value = '1.76'
$c.Style = 'Normal'
$c = $ws.Cells.Item(31, 5)
$c.NumberFormat = '@'
$c.Value = '  +5.64%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(32, 5)
$c.NumberFormat = '@'
$c.Value = '  +5.98%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = '@'
$c.Value = '6.42'
$c.Style = 'Normal'
$c = $ws.Cells.Item(33, 5)
$c.NumberFormat = '@'
$c.Value = '  +3.13%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(34, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.00%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = '@'
$c.Value = '0.994'
$c.Style = 'Normal'
$c = $ws.Cells.Item(35, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.25%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = '@'
$c.Value = '18.31'
$c.Style = 'Normal'
$c = $ws.Cells.Item(36, 5)
$c.NumberFormat = '@'
$c.Value = '  +3.54%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(37, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.34%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = '@'
$c.Value = '  +3.74%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = '@'
$c.Value = '  +5.34%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = '@'
$c.Value = '  +1.88%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = '@'
$c.Value = '0.793'
$c.Style = 'Normal'
$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = '@'
$c.Value = '  +2.43%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = '@'
$c.Value = '284.00'
$c.Style = 'Normal'
$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = '@'
$c.Value = '  +6.14%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = '@'
$c.Value = '3.52'
$c.Style = 'Normal'
$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = '@'
$c.Value = '  +4.78%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = '@'
$c.Value = '5.15'
$c.Style = 'Normal'
$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = '@'
$c.Value = '  +5.50%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = '@'
$c.Value = '133.15'
$c.Style = 'Normal'
$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = '@'
$c.Value = '  +11.16%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = '@'
$c.Value = '  +2.81%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = '@'
$c.Value = '0.0931'
$c.Style = 'Normal'
$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = '@'
$c.Value = '  +2.96%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = '@'
$c.Value = '  +6.34%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = '@'
$c.Value = '  +5.46%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = '@'
$c.Value = '  +5.14%  '
$c.Style = 'Normal'
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = '@'
$c.Value = '1.769.04'
$c.Style = 'Normal'
$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = '@'
$c.Value = '  +4.23%  '
$c.Style = 'Normal'
